$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (B2 gets a new value, C2 changes)
$ws.Range("A2").Value = "IP_Long An"
$ws.Range("B2").Value = "slp (sea logistics partners)"
$ws.Range("C2").Value = "Developer"

# Remove row 3 entirely (shifts cells up / deletes the row)
$ws.Range("A3:C3").Delete()
